$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new BOM row for the e-Ink Display (row 20; rows 17-19 left blank
# for upcoming footprint entries per the commit message)
$ws.Range("A20").Value = "e-Ink Display"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "https://www.adafruit.com/product/4777"

# Move the active selection to where the author ended up after editing
[void]$ws.Range("B25").Select()
